$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 417, pushing existing rows 417..472 down to 418..473
$ws.Rows.Item(417).Insert()

# Populate the newly inserted row 417 with its data
$ws.Cells.Item(417, 1).Value = 11
$ws.Cells.Item(417, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(417, 3).Value = "Bíobío"
$ws.Cells.Item(417, 4).Value = 45209
$ws.Cells.Item(417, 5).Value = 8
$ws.Cells.Item(417, 6).Value = 100112009
$ws.Cells.Item(417, 7).Value = "Acelga"
$ws.Cells.Item(417, 8).Value = "Sin especificar"
$ws.Cells.Item(417, 9).Value = "Primera"
$ws.Cells.Item(417, 10).Value = 220
$ws.Cells.Item(417, 11).Value = 550
$ws.Cells.Item(417, 12).Value = 600
$ws.Cells.Item(417, 13).Value = 577
$ws.Cells.Item(417, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(417, 15).Value = "Región de Ñuble"
$ws.Cells.Item(417, 16).Value = 577
$ws.Cells.Item(417, 17).Value = 1
$ws.Cells.Item(417, 18).Value = "Hortaliza"
